# This script reproduces the "added react with js" diff, which:
#   - re-styles/resizes existing rows 103-106
#   - appends new rows 107-117 with tracked coding-practice data
#   - grows the sheet dimension from A1:J106 to A1:J117
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Named constants used below (iron_native does not pre-define the usual Excel enums):
#   xlPasteFormats = -4122   (paste cell formatting only)
#   xlPasteValues  = -4163   (paste literal values only, no formulas)
$xlPasteFormats = -4122
$xlPasteValues  = -4163

# --- Step 1: rows 104 and 105 currently carry two distinctive style patterns that
#     are about to be overwritten (see Step 2 below). Clone those patterns onto the
#     new rows 109-110 and 111-113 *before* rows 104/105 themselves get re-styled. ---
$ws.Range("A104:J104").Copy()
$ws.Range("A109:J109").PasteSpecial($xlPasteFormats)
$ws.Range("A110:J110").PasteSpecial($xlPasteFormats)

$ws.Range("A105:J105").Copy()
$ws.Range("A111:J111").PasteSpecial($xlPasteFormats)
$ws.Range("A112:J112").PasteSpecial($xlPasteFormats)
$ws.Range("A113:J113").PasteSpecial($xlPasteFormats)

# --- Step 2: rows 104, 105, 106, 107 and 108 all end up sharing the style pattern
#     already present on row 103, so clone it across all of them. ---
$ws.Range("A103:J103").Copy()
$ws.Range("A104:J104").PasteSpecial($xlPasteFormats)
$ws.Range("A105:J105").PasteSpecial($xlPasteFormats)
$ws.Range("A106:J106").PasteSpecial($xlPasteFormats)
$ws.Range("A107:J107").PasteSpecial($xlPasteFormats)
$ws.Range("A108:J108").PasteSpecial($xlPasteFormats)

# --- Step 3: fill in the cell values/text for every touched row (106 through 117).
#     Text-like values (dates, clock times, and the odd numeric-looking-but-text Wpm
#     entries in rows 114-117) are entered as a quoted formula and immediately frozen
#     into a literal value in-place; this stops Excel from auto-converting date- or
#     time-shaped text into real date/time serial numbers. ---
# row 106
$c = $ws.Range("A106"); $c.Formula = "=""2025-02-26"""; $c.Copy(); $c.PasteSpecial($xlPasteValues)
$c = $ws.Range("B106"); $c.Formula = "=""8:30"""; $c.Copy(); $c.PasteSpecial($xlPasteValues)
$ws.Range("C106").Value = 35
$c = $ws.Range("D106"); $c.Formula = "=""4:51"""; $c.Copy(); $c.PasteSpecial($xlPasteValues)
$c = $ws.Range("E106"); $c.Formula = "=""3:41"""; $c.Copy(); $c.PasteSpecial($xlPasteValues)
$ws.Range("F106").Value = 0
$ws.Range("G106").Value = 0
$ws.Range("H106").Value = 209
$ws.Range("I106").Value = 209
$ws.Range("J106").Value = 166

# row 107
$c = $ws.Range("A107"); $c.Formula = "=""2025-02-27"""; $c.Copy(); $c.PasteSpecial($xlPasteValues)
$c = $ws.Range("B107"); $c.Formula = "=""8:30"""; $c.Copy(); $c.PasteSpecial($xlPasteValues)
$ws.Range("C107").Value = 40
$c = $ws.Range("D107"); $c.Formula = "=""7:45"""; $c.Copy(); $c.PasteSpecial($xlPasteValues)
$c = $ws.Range("E107"); $c.Formula = "=""5:37"""; $c.Copy(); $c.PasteSpecial($xlPasteValues)
$ws.Range("F107").Value = 0
$ws.Range("G107").Value = 8
$ws.Range("H107").Value = 188
$ws.Range("I107").Value = 196
$ws.Range("J107").Value = 167

# row 108
$c = $ws.Range("A108"); $c.Formula = "=""2025-02-28"""; $c.Copy(); $c.PasteSpecial($xlPasteValues)
$c = $ws.Range("B108"); $c.Formula = "=""10:20"""; $c.Copy(); $c.PasteSpecial($xlPasteValues)
$ws.Range("C108").Value = 40
$c = $ws.Range("D108"); $c.Formula = "=""4:20"""; $c.Copy(); $c.PasteSpecial($xlPasteValues)
$c = $ws.Range("E108"); $c.Formula = "=""2:06"""; $c.Copy(); $c.PasteSpecial($xlPasteValues)
$ws.Range("F108").Value = 0
$ws.Range("G108").Value = 4
$ws.Range("H108").Value = 140
$ws.Range("I108").Value = 144
$ws.Range("J108").Value = 168

# row 109
$c = $ws.Range("A109"); $c.Formula = "=""2025-03-01"""; $c.Copy(); $c.PasteSpecial($xlPasteValues)
$c = $ws.Range("B109"); $c.Formula = "=""9:00"""; $c.Copy(); $c.PasteSpecial($xlPasteValues)
$ws.Range("C109").Value = 37
$c = $ws.Range("D109"); $c.Formula = "=""1:35"""; $c.Copy(); $c.PasteSpecial($xlPasteValues)
$c = $ws.Range("E109"); $c.Formula = "=""0:17"""; $c.Copy(); $c.PasteSpecial($xlPasteValues)
$ws.Range("F109").Value = 0
$ws.Range("G109").Value = 0
$ws.Range("H109").Value = 19
$ws.Range("I109").Value = 19
$ws.Range("J109").Value = 169

# row 110
$c = $ws.Range("A110"); $c.Formula = "=""2025-03-03"""; $c.Copy(); $c.PasteSpecial($xlPasteValues)
$c = $ws.Range("B110"); $c.Formula = "=""8:41"""; $c.Copy(); $c.PasteSpecial($xlPasteValues)
$ws.Range("C110").Value = 39
$c = $ws.Range("D110"); $c.Formula = "=""3:41"""; $c.Copy(); $c.PasteSpecial($xlPasteValues)
$c = $ws.Range("E110"); $c.Formula = "=""0:56"""; $c.Copy(); $c.PasteSpecial($xlPasteValues)
$ws.Range("F110").Value = 0
$ws.Range("G110").Value = 0
$ws.Range("H110").Value = 80
$ws.Range("I110").Value = 80
$ws.Range("J110").Value = 170

# row 111
$c = $ws.Range("A111"); $c.Formula = "=""2025-03-04"""; $c.Copy(); $c.PasteSpecial($xlPasteValues)
$c = $ws.Range("B111"); $c.Formula = "=""9:00"""; $c.Copy(); $c.PasteSpecial($xlPasteValues)
$ws.Range("C111").Value = 41
$c = $ws.Range("D111"); $c.Formula = "=""3:35"""; $c.Copy(); $c.PasteSpecial($xlPasteValues)
$c = $ws.Range("E111"); $c.Formula = "=""0:52"""; $c.Copy(); $c.PasteSpecial($xlPasteValues)
$ws.Range("F111").Value = 0
$ws.Range("G111").Value = 0
$ws.Range("H111").Value = 58
$ws.Range("I111").Value = 58
$ws.Range("J111").Value = 171

# row 112
$c = $ws.Range("A112"); $c.Formula = "=""2025-03-05"""; $c.Copy(); $c.PasteSpecial($xlPasteValues)
$c = $ws.Range("B112"); $c.Formula = "=""8:30"""; $c.Copy(); $c.PasteSpecial($xlPasteValues)
$ws.Range("C112").Value = 44
$c = $ws.Range("D112"); $c.Formula = "=""5:24"""; $c.Copy(); $c.PasteSpecial($xlPasteValues)
$c = $ws.Range("E112"); $c.Formula = "=""2:36"""; $c.Copy(); $c.PasteSpecial($xlPasteValues)
$ws.Range("F112").Value = 0
$ws.Range("G112").Value = 0
$ws.Range("H112").Value = 200
$ws.Range("I112").Value = 200
$ws.Range("J112").Value = 172

# row 113
$c = $ws.Range("A113"); $c.Formula = "=""2025-03-06"""; $c.Copy(); $c.PasteSpecial($xlPasteValues)
$c = $ws.Range("B113"); $c.Formula = "=""8:30"""; $c.Copy(); $c.PasteSpecial($xlPasteValues)
$ws.Range("C113").Value = 42
$c = $ws.Range("D113"); $c.Formula = "=""5:24"""; $c.Copy(); $c.PasteSpecial($xlPasteValues)
$c = $ws.Range("E113"); $c.Formula = "=""2:46"""; $c.Copy(); $c.PasteSpecial($xlPasteValues)
$ws.Range("F113").Value = 0
$ws.Range("G113").Value = 0
$ws.Range("H113").Value = 208
$ws.Range("I113").Value = 208
$ws.Range("J113").Value = 173

# row 114
$c = $ws.Range("A114"); $c.Formula = "=""2025-03-07"""; $c.Copy(); $c.PasteSpecial($xlPasteValues)
$c = $ws.Range("B114"); $c.Formula = "=""6:58"""; $c.Copy(); $c.PasteSpecial($xlPasteValues)
$c = $ws.Range("C114"); $c.Formula = "=""40"""; $c.Copy(); $c.PasteSpecial($xlPasteValues)
$c = $ws.Range("D114"); $c.Formula = "=""5:58"""; $c.Copy(); $c.PasteSpecial($xlPasteValues)
$c = $ws.Range("E114"); $c.Formula = "=""3:42"""; $c.Copy(); $c.PasteSpecial($xlPasteValues)
$ws.Range("F114").Value = 0
$ws.Range("G114").Value = 0
$ws.Range("H114").Value = 113
$ws.Range("I114").Value = 113
$ws.Range("J114").Value = 174

# row 115
$c = $ws.Range("A115"); $c.Formula = "=""2025-03-09"""; $c.Copy(); $c.PasteSpecial($xlPasteValues)
$c = $ws.Range("B115"); $c.Formula = "=""7:27"""; $c.Copy(); $c.PasteSpecial($xlPasteValues)
$c = $ws.Range("C115"); $c.Formula = "=""40"""; $c.Copy(); $c.PasteSpecial($xlPasteValues)
$c = $ws.Range("D115"); $c.Formula = "=""6:27"""; $c.Copy(); $c.PasteSpecial($xlPasteValues)
$c = $ws.Range("E115"); $c.Formula = "=""0:14"""; $c.Copy(); $c.PasteSpecial($xlPasteValues)
$ws.Range("F115").Value = 0
$ws.Range("G115").Value = 0
$ws.Range("H115").Value = 0
$ws.Range("I115").Value = 0
$ws.Range("J115").Value = 175

# row 116
$c = $ws.Range("A116"); $c.Formula = "=""2025-03-10"""; $c.Copy(); $c.PasteSpecial($xlPasteValues)
$c = $ws.Range("B116"); $c.Formula = "=""8:30"""; $c.Copy(); $c.PasteSpecial($xlPasteValues)
$c = $ws.Range("C116"); $c.Formula = "=""40"""; $c.Copy(); $c.PasteSpecial($xlPasteValues)
$c = $ws.Range("D116"); $c.Formula = "=""6:39"""; $c.Copy(); $c.PasteSpecial($xlPasteValues)
$c = $ws.Range("E116"); $c.Formula = "=""4:15"""; $c.Copy(); $c.PasteSpecial($xlPasteValues)
$ws.Range("F116").Value = 0
$ws.Range("G116").Value = 0
$ws.Range("H116").Value = 228
$ws.Range("I116").Value = 228
$ws.Range("J116").Value = 176

# row 117
$c = $ws.Range("A117"); $c.Formula = "=""2025-03-11"""; $c.Copy(); $c.PasteSpecial($xlPasteValues)
$c = $ws.Range("B117"); $c.Formula = "=""8:00"""; $c.Copy(); $c.PasteSpecial($xlPasteValues)
$c = $ws.Range("C117"); $c.Formula = "=""43"""; $c.Copy(); $c.PasteSpecial($xlPasteValues)
$c = $ws.Range("D117"); $c.Formula = "=""5:22"""; $c.Copy(); $c.PasteSpecial($xlPasteValues)
$c = $ws.Range("E117"); $c.Formula = "=""2:32"""; $c.Copy(); $c.PasteSpecial($xlPasteValues)
$ws.Range("F117").Value = 1
$ws.Range("G117").Value = 0
$ws.Range("H117").Value = 164
$ws.Range("I117").Value = 165
$ws.Range("J117").Value = 177

# --- Step 4: row heights (19.5pt for rows 103-106, 16.5pt for the newly added
#     rows 107-113; rows 114-117 keep the default height, just like the diff) ---
$ws.Rows.Item(103).RowHeight = 19.5
$ws.Rows.Item(104).RowHeight = 19.5
$ws.Rows.Item(105).RowHeight = 19.5
$ws.Rows.Item(106).RowHeight = 19.5
$ws.Rows.Item(107).RowHeight = 16.5
$ws.Rows.Item(108).RowHeight = 16.5
$ws.Rows.Item(109).RowHeight = 16.5
$ws.Rows.Item(110).RowHeight = 16.5
$ws.Rows.Item(111).RowHeight = 16.5
$ws.Rows.Item(112).RowHeight = 16.5
$ws.Rows.Item(113).RowHeight = 16.5

